# Updates CATALOGO_GRUPOS.xlsx:
#  - "MEJORAR" sheet (sheet1): inserts several new product codes at the top
#    and at the end of the list (grows from 10 to 40 rows), and centers the
#    header cell (A1) by adding center alignment on top of its existing font.
#  - "PREMIUM" sheet (sheet2): fixes one product code typo
#    (evorieg153 -> evorieg0153); all other values stay the same.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "MEJORAR"
# ---------------------------------------------------------------------
$wsMejorar = $wb.Worksheets.Item("MEJORAR")

$mejorarValues = @(
    "Codigo_Producto",
    "SIM18310",
    "KIMERAC1",
    "PROBOT41",
    "evorieg0153",
    "EGWX 01",
    "EGWX 02",
    "PX120314",
    "evol0088",
    "evol1000",
    "evol0330",
    "PERFA0261",
    "GAG12103AR",
    "TOR01522",
    "evol0088",
    "evol0330",
    "evol3089",
    "evol0025",
    "evol3245",
    "evol1970",
    "evo115la",
    "evo115tu",
    "evol0028",
    "evol3510",
    "evol0070",
    "evol2530",
    "evol0107",
    "evol0435",
    "evol0111",
    "evol0088",
    "evol3970",
    "evol0177",
    "evol0174",
    "evol2205",
    "evol1361",
    "evol3210",
    "evol5530",
    "evo115co",
    "CON205",
    "TF414"
)

for ($i = 0; $i -lt $mejorarValues.Length; $i++) {
    $wsMejorar.Cells.Item($i + 1, 1).Value = $mejorarValues[$i]
}

# Header keeps its Consolas font but gains centered horizontal alignment.
$wsMejorar.Cells.Item(1, 1).HorizontalAlignment = -4108

# All data rows (2..40) are center-aligned, matching the rest of the
# catalogue - re-apply explicitly since the newly added rows (11..40)
# don't inherit the formatting of the original 10-row range.
$wsMejorar.Range("A2:A40").HorizontalAlignment = -4108

# Matches the author's final selection/scroll state on this sheet.
$wsMejorar.Activate() | Out-Null
$wsMejorar.Range("A22").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "PREMIUM"
# ---------------------------------------------------------------------
$wsPremium = $wb.Worksheets.Item("PREMIUM")

$premiumValues = @(
    "Codigo_Producto",
    "evol1000",
    "evol0330",
    "PERFA0261",
    "evol0088",
    "evo115tu",
    "evol0028",
    "evol3510",
    "evorieg0153",
    "evol0070",
    "evol2530",
    "evol0107",
    "evol0435"
)

for ($i = 0; $i -lt $premiumValues.Length; $i++) {
    $wsPremium.Cells.Item($i + 1, 1).Value = $premiumValues[$i]
}

$wsPremium.Activate() | Out-Null
$wsPremium.Range("A10").Select() | Out-Null

# Re-activate the sheet that was selected in the saved file.
$wsMejorar.Activate() | Out-Null
